$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    $ws.Range("B2").Value = 0.01384957621052772
    $ws.Range("C2").Value = 2.21878651858606
    $ws.Range("D2").Value = 21.84152086182795
    $ws.Range("E2").Value = 4.673491292580734
    $ws.Range("F2").Value = 4.783449836594167
    $ws.Range("G2").Value = 22
    $ws.Range("B3").Value = -0.1108950891496856
    $ws.Range("C3").Value = 2.52800863553415
    $ws.Range("D3").Value = 20.81929224357433
    $ws.Range("E3").Value = 4.562816262307122
    $ws.Range("F3").Value = 4.674114274267978
    $ws.Range("G3").Value = 21
    $ws.Range("B4").Value = -0.5909552399615555
    $ws.Range("C4").Value = 1.803240045667984
    $ws.Range("D4").Value = 9.880028153484647
    $ws.Range("E4").Value = 3.143251207505478
    $ws.Range("F4").Value = 3.167399740090089
    $ws.Range("G4").Value = 20
    $ws.Range("B5").Value = -0.1812638745297453
    $ws.Range("C5").Value = 1.430346199258623
    $ws.Range("D5").Value = 9.695185292776637
    $ws.Range("E5").Value = 3.113709249878132
    $ws.Range("F5").Value = 3.193606854246072
    $ws.Range("G5").Value = 19
    $ws.Range("B6").Value = -0.1687010319474557
    $ws.Range("C6").Value = 1.931467188986993
    $ws.Range("D6").Value = 13.19207325907522
    $ws.Range("E6").Value = 3.632089379279539
    $ws.Range("F6").Value = 3.733355516202483
    $ws.Range("G6").Value = 18
    $ws.Range("B7").Value = -0.2974795460340417
    $ws.Range("C7").Value = 1.914100224451346
    $ws.Range("D7").Value = 11.03882338163619
    $ws.Range("E7").Value = 3.322472480192453
    $ws.Range("F7").Value = 3.410971252101159
    $ws.Range("G7").Value = 17
    $ws.Range("B8").Value = -0.1960622915985716
    $ws.Range("C8").Value = 1.836194656235551
    $ws.Range("D8").Value = 11.49912085291426
    $ws.Range("E8").Value = 3.391035365919127
    $ws.Range("F8").Value = 3.496387534695756
    $ws.Range("G8").Value = 16
    $ws.Range("B9").Value = -0.08754086917198256
    $ws.Range("C9").Value = 1.885964882007152
    $ws.Range("D9").Value = 11.84300447804472
    $ws.Range("E9").Value = 3.441366658472288
    $ws.Range("F9").Value = 3.561000221788015
    $ws.Range("G9").Value = 15
    $ws.Range("B10").Value = -0.1199187714296267
    $ws.Range("C10").Value = 2.158918608588491
    $ws.Range("D10").Value = 14.23208708490037
    $ws.Range("E10").Value = 3.772543847975842
    $ws.Range("F10").Value = 3.912975378092231
    $ws.Range("G10").Value = 14
    $ws.Range("B11").Value = -0.08524854279770107
    $ws.Range("C11").Value = 2.086633922413215
    $ws.Range("D11").Value = 13.79239774298984
    $ws.Range("E11").Value = 3.713811753843999
    $ws.Range("F11").Value = 3.864439325182775
    $ws.Range("G11").Value = 13
